# "corregido error en base de datos" - fix the match results stored on the
# "datos" sheet (rows 82-91) and update the active sheet/selection state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("datos")

# Corrected match results for rows 82-91 (equipo1, goles1, equipo2, goles2)
$ws1.Range("B82").Value = "Arsenal"
$ws1.Range("C82").Value = 3
$ws1.Range("D82").Value = "Godoy Cruz Mza."
$ws1.Range("E82").Value = 0

$ws1.Range("B83").Value = "At. Rafaela"
$ws1.Range("C83").Value = 0
$ws1.Range("D83").Value = "Olimpo (BB)"
$ws1.Range("E83").Value = 0

$ws1.Range("B84").Value = "Belgrano (Cba)"
$ws1.Range("C84").Value = 1
$ws1.Range("D84").Value = "Racing Club"
$ws1.Range("E84").Value = 4

$ws1.Range("B85").Value = "Boca Juniors"
$ws1.Range("C85").Value = 1
$ws1.Range("D85").Value = "Quilmes"
$ws1.Range("E85").Value = 0

$ws1.Range("B86").Value = "Estudiantes LP"
$ws1.Range("C86").Value = 3
$ws1.Range("D86").Value = "Vélez Sarsfield"
$ws1.Range("E86").Value = 2

$ws1.Range("B87").Value = "Independiente"
$ws1.Range("C87").Value = 2
$ws1.Range("D87").Value = "Rosario Central"
$ws1.Range("E87").Value = 0

$ws1.Range("B88").Value = "Lanús"
$ws1.Range("C88").Value = 1
$ws1.Range("D88").Value = "River Plate"
$ws1.Range("E88").Value = 1

$ws1.Range("B89").Value = "Newell's"
$ws1.Range("C89").Value = 0
$ws1.Range("D89").Value = "Banfield"
$ws1.Range("E89").Value = 3

$ws1.Range("B90").Value = "San Lorenzo"
$ws1.Range("C90").Value = 0
$ws1.Range("D90").Value = "Gimnasia LP"
$ws1.Range("E90").Value = 2

$ws1.Range("B91").Value = "Tigre"
$ws1.Range("C91").Value = 2
$ws1.Range("D91").Value = "Def. y Justicia"
$ws1.Range("E91").Value = 1

# Make "datos" the active sheet/tab and move the selection to F3
# (this also clears the previously active "equipos" tab's tabSelected flag).
$ws1.Activate()
$ws1.Range("F3").Select()
